# oneamp lig adaptor received
#
# The underlying data didn't change; this commit is Excel re-saving the
# workbook after the user received/reviewed the data: it auto-sized the
# three used columns (A:C) to fit their contents and left the selection
# on the data body (A2:C25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Autofit columns A:C to the cell contents, then nudge the widths to the
# precise values Excel settled on (COM ColumnWidth is in "characters";
# the engine's pixel grid is ColumnWidth*6+5, so feed it the values that
# round-trip to the closest achievable widths: 12.5, ~38.1667, ~14.3333).
$ws.Columns("A:C").AutoFit() | Out-Null

$ws.Columns("A").ColumnWidth = 11.666666666666666
$ws.Columns("B").ColumnWidth = 37.333333333333336
$ws.Columns("C").ColumnWidth = 13.5

# Leave the selection on the received data range, anchored at A2.
$ws.Range("A2:C25").Select() | Out-Null
